# Insert a new weekly price record at row 11 (pushing existing rows 11-62
# down to 12-63, which also grows the sheet's used range to A1:R63).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(11).Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(11, 1).Value = 1
$ws.Cells.Item(11, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(11, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(11, 4).Value = 44649
$ws.Cells.Item(11, 5).Value = 15
$ws.Cells.Item(11, 6).Value = 100112012
$ws.Cells.Item(11, 7).Value = "Espinaca"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 300
$ws.Cells.Item(11, 11).Value = 1800
$ws.Cells.Item(11, 12).Value = 2000
$ws.Cells.Item(11, 13).Value = 1900
$ws.Cells.Item(11, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(11, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(11, 16).Value = 633
$ws.Cells.Item(11, 17).Value = 3
$ws.Cells.Item(11, 18).Value = "Hortaliza"
